$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns remain plain text so numeric-looking strings
# (e.g. "1.013", "0.00001108") are not auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.612.48"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.009.08"
$ws.Range("E3").Value = "  -5.09%  "
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "330.79"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "0.4969"
$ws.Range("E7").Value = "  -4.38%  "
$ws.Range("D8").Value = "0.4224"
$ws.Range("E8").Value = "  -4.84%  "
$ws.Range("D9").Value = "53.77"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").Value = "0.08851"
$ws.Range("E10").Value = "  -5.13%  "
$ws.Range("D11").Value = "1.124"
$ws.Range("E11").Value = "  -4.88%  "
$ws.Range("D12").Value = "2.169.91"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").Value = "23.05"
$ws.Range("E13").Value = "  -8.36%  "
$ws.Range("D14").Value = "8.146"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").Value = "6.518"
$ws.Range("E15").Value = "  -5.51%  "
$ws.Range("D16").Value = "96.65"
$ws.Range("E16").Value = "  -6.26%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.00001108"
$ws.Range("E18").Value = "  -4.59%  "
$ws.Range("D19").Value = "0.06622"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "19.59"
$ws.Range("E20").Value = "  -8.84%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "6.012"
$ws.Range("E22").Value = "  -4.60%  "
$ws.Range("D23").Value = "29.623.20"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("E24").Value = "  -6.02%  "
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D26").Value = "2.218.75"
$ws.Range("E26").Value = "  -5.56%  "
$ws.Range("D27").Value = "158.11"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("D29").Value = "6.569"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "2.340"
$ws.Range("E30").Value = "  -7.81%  "
$ws.Range("D31").Value = "127.37"
$ws.Range("E31").Value = "  -5.13%  "
$ws.Range("D32").Value = "1.059"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("D33").Value = "0.09975"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("D34").Value = "1.557"
$ws.Range("E34").Value = "  -11.81%  "
$ws.Range("D35").Value = "3.809"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("E36").Value = "  -6.97%  "
$ws.Range("D37").Value = "9.623"
$ws.Range("E37").Value = "  -10.42%  "
$ws.Range("D38").Value = "0.02464"
$ws.Range("E38").Value = "  -6.11%  "
$ws.Range("D39").Value = "0.06406"
$ws.Range("E39").Value = "  -6.58%  "
$ws.Range("D40").Value = "1.295"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "11.85"
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("D42").Value = "0.6528"
$ws.Range("E42").Value = "  -7.75%  "
$ws.Range("D43").Value = "0.2077"
$ws.Range("E43").Value = "  -7.45%  "
$ws.Range("D44").Value = "1.011"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "0.6346"
$ws.Range("E45").Value = "  -7.27%  "
$ws.Range("D46").Value = "2.233"
$ws.Range("E46").Value = "  -5.23%  "
$ws.Range("D47").Value = "13.53"
$ws.Range("E47").Value = "  -6.82%  "
$ws.Range("D48").Value = "1.269"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "3.562"
$ws.Range("E49").Value = "  -2.11%  "

# Rows 50 and 51 swap coin order (BabyDogeCoin now ranks ahead of Cronos)
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.00000000328"
$ws.Range("E50").Value = "  -8.36%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.07027"
$ws.Range("E51").Value = "  -1.22%  "
